# Weekly price-sheet update: insert a new daily record for Kiwi
# (Vega Monumental Concepción) as row 165, pushing the existing
# rows 165-198 down to 166-199.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 165 - shifts rows 165:198 down to 166:199
# and extends the used range to row 199 automatically.
$ws.Rows(165).Insert()

# Populate the newly inserted row 165 with the new price record.
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 44785
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100101
$ws.Cells.Item(165, 8).Value = "Berries"
$ws.Cells.Item(165, 9).Value = 100101007
$ws.Cells.Item(165, 10).Value = "Kiwi"
$ws.Cells.Item(165, 11).Value = "Hayward"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 220
$ws.Cells.Item(165, 14).Value = 6000
$ws.Cells.Item(165, 15).Value = 6500
$ws.Cells.Item(165, 16).Value = 6273
$ws.Cells.Item(165, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 348
$ws.Cells.Item(165, 20).Value = 18

# Make sure the date cell keeps the workbook's date number format, same
# as every other cell in column D.
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
